# Workbook was regenerated with an older Excel build: column "Danish term"
# (B) was removed, the "(system)API Search Form" / "(system)API Search Form
# Sorting" columns were replaced by a single "(system)API Property Mapping"
# column, and the header row was reselected/resized.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the header row left by one (the old "Danish term" column is gone)
# and fold the two old "(system)API Search Form..." columns into a single
# new "(system)API Property Mapping" column placed right before
# "(system)API Search Criteria Mapping".
$ws.Range("B1").Value = "Imported from CMS (Y/N)"
$ws.Range("C1").Value = "Optionality  (O/M)CO is imported from CMS"
$ws.Range("D1").Value = "Optionality  (O/M)CO is created in CS"
$ws.Range("E1").Value = "Editing  (E/R/H)CO is imported from CMS"
$ws.Range("F1").Value = "Editing  (E/R/H)CO is created in CS"
$ws.Range("G1").Value = "Data Type"
$ws.Range("H1").Value = "Controlled  vocabulary"
$ws.Range("I1").Value = "Searchable"
$ws.Range("J1").Value = "Description"
$ws.Range("K1").Value = "Sample Data"
$ws.Range("L1").Value = "Comment"
$ws.Range("M1").Value = "(system)CS URI"
$ws.Range("N1").Value = "(system)CS Definition Property Mapping"
$ws.Range("O1").Value = "(system)API Property Mapping"
$ws.Range("P1").Value = "(system)API Search Criteria Mapping"
$ws.Range("Q1").ClearContents()
$ws.Range("R1").ClearContents()

# Header row now renders slightly shorter, and the whole row is selected
# (matches the resave from the older Excel build captured in the diff).
$ws.Rows(1).RowHeight = 14.25
$null = $ws.Rows(1).Select()
